$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").ClearContents()
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = 0
$ws.Range("H76").Value = 4360.885
$ws.Range("I76").Value = 14249.5
$ws.Range("J76").Value = 3536.8333
$ws.Range("K76").Value = 14249.5
$ws.Range("L76").Value = 3536.8333
$ws.Range("M76").Value = -13934.5
$ws.Range("N76").Value = -4166.8333
$ws.Range("H79").Value = 4360.885
$ws.Range("I79").Value = 14249.5
$ws.Range("J79").Value = 3536.8333
$ws.Range("K79").Value = 14249.5
$ws.Range("L79").Value = 3536.8333
$ws.Range("M79").Value = -13157.5
$ws.Range("N79").Value = -5720.8333
$ws.Range("H135").Value = 964.61365
$ws.Range("I135").Value = 869.4167
$ws.Range("J135").Value = 1393
$ws.Range("K135").Value = 7824.7503
$ws.Range("L135").Value = 12537
$ws.Range("M135").Value = -5289.7503
$ws.Range("N135").Value = -17607
$ws.Range("H138").Value = 3925.4
$ws.Range("I138").Value = 2549.375
$ws.Range("J138").Value = 4089.7014
$ws.Range("K138").Value = 7648.125
$ws.Range("L138").Value = 12269.1042
$ws.Range("M138").Value = -2508.125
$ws.Range("N138").Value = -22549.1042

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H5").Value = 49
$ws.Range("I5").Value = 49
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 49
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 63
$ws.Range("H32").Value = 25670090
$ws.Range("I32").Value = 66681068
$ws.Range("J32").Value = 38229.168
$ws.Range("K32").Value = 66681068
$ws.Range("L32").Value = 38229.168
$ws.Range("M32").Value = -66680781
$ws.Range("N32").Value = -38803.168
$ws.Range("H45").Value = 1459.5
$ws.Range("I45").Value = 815.4286
$ws.Range("J45").Value = 1596.1212
$ws.Range("K45").Value = 815.4286
$ws.Range("L45").Value = 1596.1212
$ws.Range("M45").Value = -438.4286
$ws.Range("N45").Value = -2350.1212
$ws.Range("H61").Value = 2641.889
$ws.Range("I61").Value = 1391.7778
$ws.Range("J61").Value = 3892
$ws.Range("K61").Value = 1391.7778
$ws.Range("L61").Value = 3892
$ws.Range("M61").Value = -1179.7778
$ws.Range("N61").Value = -4316
$ws.Range("H122").Value = 1399.2142
$ws.Range("I122").Value = 1194.9166
$ws.Range("J122").Value = 2625
$ws.Range("K122").Value = 3584.7498
$ws.Range("L122").Value = 7875
$ws.Range("M122").Value = -1134.7498
$ws.Range("N122").Value = -12775
$ws.Range("H136").Value = 2641.889
$ws.Range("I136").Value = 1391.7778
$ws.Range("J136").Value = 3892
$ws.Range("K136").Value = 4175.3334
$ws.Range("L136").Value = 11676
$ws.Range("M136").Value = -1625.3334
$ws.Range("N136").Value = -16776

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H4").Value = 49
$ws.Range("I4").Value = 49
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 49
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 66
$ws.Range("H107").Value = 1741.6666
$ws.Range("I107").Value = 1650
$ws.Range("J107").Value = 1925
$ws.Range("K107").Value = 1650
$ws.Range("L107").Value = 1925
$ws.Range("M107").Value = 270
$ws.Range("N107").Value = -5765

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H94").Value = 3148.7188
$ws.Range("I94").Value = 1549.9166
$ws.Range("J94").Value = 4108
$ws.Range("K94").Value = 1549.9166
$ws.Range("L94").Value = 4108
$ws.Range("M94").Value = -1098.9166
$ws.Range("N94").Value = -5010
$ws.Range("H123").Value = 46833.332
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 46833.332
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 46833.332
$ws.Range("N123").Value = -56633.332
$ws.Range("H129").Value = 28695
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 28695
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 28695
$ws.Range("N129").Value = -38695
$ws.Range("H141").Value = 43075.062
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 43075.062
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 43075.062
$ws.Range("N141").Value = -53435.062

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H80").Value = 2430.1
$ws.Range("I80").Value = 1702
$ws.Range("J80").Value = 2511
$ws.Range("K80").Value = 5106
$ws.Range("L80").Value = 7533
$ws.Range("M80").Value = -4170
$ws.Range("N80").Value = -9405
$ws.Range("H83").Value = 2430.1
$ws.Range("I83").Value = 1702
$ws.Range("J83").Value = 2511
$ws.Range("K83").Value = 15318
$ws.Range("L83").Value = 22599
$ws.Range("M83").Value = -10638
$ws.Range("N83").Value = -31959
$ws.Range("H131").Value = 619.86664
$ws.Range("I131").Value = 363.75
$ws.Range("J131").Value = 912.5714
$ws.Range("K131").Value = 1091.25
$ws.Range("L131").Value = 2737.7142
$ws.Range("M131").Value = 3948.75
$ws.Range("N131").Value = -12817.7142

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H68").Value = 2135.0435
$ws.Range("I68").Value = 1988.5883
$ws.Range("J68").Value = 2550
$ws.Range("K68").Value = 1988.5883
$ws.Range("L68").Value = 2550
$ws.Range("M68").Value = -1239.5883
$ws.Range("N68").Value = -4048
$ws.Range("H71").Value = 2135.0435
$ws.Range("I71").Value = 1988.5883
$ws.Range("J71").Value = 2550
$ws.Range("K71").Value = 9942.941499999999
$ws.Range("L71").Value = 12750
$ws.Range("M71").Value = -6198.941499999999
$ws.Range("N71").Value = -20238
$ws.Range("H74").Value = 26750
$ws.Range("I74").Value = 26750
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 26750
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -25752
$ws.Range("H77").Value = 26750
$ws.Range("I77").Value = 26750
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 80250
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -75258

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").ClearContents()
$ws.Range("N64").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").ClearContents()
$ws.Range("N67").Value = 0
$ws.Range("H81").Value = 2546.6667
$ws.Range("I81").Value = 2546.6667
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 5093.3334
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -4032.3334
$ws.Range("H84").Value = 2546.6667
$ws.Range("I84").Value = 2546.6667
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 25466.667
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -20162.667
$ws.Range("H113").Value = 527.8
$ws.Range("I113").Value = 199
$ws.Range("J113").Value = 610
$ws.Range("K113").Value = 597
$ws.Range("L113").Value = 1830
$ws.Range("M113").Value = 1573
$ws.Range("N113").Value = -6170
$ws.Range("H126").Value = 4647.5
$ws.Range("I126").Value = 4158.5
$ws.Range("J126").Value = 5625.5
$ws.Range("K126").Value = 12475.5
$ws.Range("L126").Value = 16876.5
$ws.Range("M126").Value = -10005.5
$ws.Range("N126").Value = -21816.5
